# Fixing network data cleaning scripts
# - rename header columns to snake_case field names
# - normalize "de/del/el/los/la" -> "De/Del/El/Los/La" in municipality/state names
# - drop the trailing metadata/footer rows (181-480), shrinking used range to A1:D179

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Capitalize lowercase Spanish prepositions in place names ---
$ws.Range("B8").Value = "Amatenango De La Frontera"
$ws.Range("B10").Value = "Comitán De Domínguez"
$ws.Range("B16").Value = "Salto De Agua"
$ws.Range("A26").Value = "Ciudad De México"
$ws.Range("A39").Value = "Estado De México"
$ws.Range("B49").Value = "Jaral Del Progreso"
$ws.Range("B53").Value = "Acapulco De Juárez"
$ws.Range("B55").Value = "Alcozauca De Guerrero"
$ws.Range("B57").Value = "Atoyac De Álvarez"
$ws.Range("B58").Value = "Chilapa De Álvarez"
$ws.Range("B59").Value = "Chilpancingo De Los Bravo"
$ws.Range("B61").Value = "Cutzamala De Pinzón"
$ws.Range("B68").Value = "Tlapa De Comonfort"
$ws.Range("B72").Value = "Atotonilco El Grande"
$ws.Range("B74").Value = "Pachuca De Soto"
$ws.Range("B78").Value = "Lagos De Moreno"
$ws.Range("B84").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B87").Value = "Tlaltizapán De Zapata"
$ws.Range("B94").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B95").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B96").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B97").Value = "Oaxaca De Juárez"
$ws.Range("B98").Value = "Ocotlán De Morelos"
$ws.Range("B99").Value = "Putla Villa De Guerrero"
$ws.Range("B109").Value = "Santa Inés Del Monte"
$ws.Range("B123").Value = "Izúcar De Matamoros"
$ws.Range("B124").Value = "Los Reyes De Juárez"
$ws.Range("B132").Value = "Tepanco De López"
$ws.Range("B133").Value = "Tetela De Ocampo"
$ws.Range("B134").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B137").Value = "Xayacatlán De Bravo"
$ws.Range("B140").Value = "Cadereyta De Montes"
$ws.Range("B160").Value = "San Pablo Del Monte"
$ws.Range("B162").Value = "Amatlán De Los Reyes"
$ws.Range("B169").Value = "Martínez De La Torre"
$ws.Range("B172").Value = "Poza Rica De Hidalgo"

# --- Drop the trailing metadata/footer block (rows 181-480) ---
$ws.Rows("181:480").Delete()
